$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 9.573375
$ws.Range("H2").Value = 28.720125
$ws.Range("I2").Value = 0.1037691388643484
$ws.Range("J2").Value = 0.1037691388643484
$ws.Range("M2").Value = 227.11144
$ws.Range("N2").Value = 681.33432
$ws.Range("O2").Value = 0.8625743548356182
$ws.Range("P2").Value = 0.8625743548356182
$ws.Range("Q2").Value = 2174.22298191
$ws.Range("R2").Value = 19568.00683719
$ws.Range("S2").Value = 0.08950859800776299
$ws.Range("T2").Value = 0.08950859800776299
# Row 3
$ws.Range("G3").Value = 9.573375
$ws.Range("H3").Value = 28.720125
$ws.Range("I3").Value = 0.1037691388643484
$ws.Range("J3").Value = 0.1037691388643484
$ws.Range("O3").Value = 0.001598666154760757
$ws.Range("P3").Value = 0.001598666154760757
$ws.Range("Q3").Value = 4.02963138725
$ws.Range("R3").Value = 36.26668248525
$ws.Range("S3").Value = 0.0001658922102111029
$ws.Range("T3").Value = 0.0001658922102111028
# Row 4
$ws.Range("G4").Value = 9.573375
$ws.Range("H4").Value = 28.720125
$ws.Range("I4").Value = 0.1037691388643484
$ws.Range("J4").Value = 0.1037691388643484
$ws.Range("M4").Value = 3.233093
$ws.Range("N4").Value = 9.699279000000001
$ws.Range("O4").Value = 0.01227935989749593
$ws.Range("P4").Value = 0.01227935989749593
$ws.Range("Q4").Value = 30.951611698875
$ws.Range("R4").Value = 278.564505289875
$ws.Range("S4").Value = 0.001274218602368566
$ws.Range("T4").Value = 0.001274218602368566
# Row 5
$ws.Range("G5").Value = 9.573375
$ws.Range("H5").Value = 28.720125
$ws.Range("I5").Value = 0.1037691388643484
$ws.Range("J5").Value = 0.1037691388643484
$ws.Range("M5").Value = 32.52945966666667
$ws.Range("N5").Value = 97.588379
$ws.Range("O5").Value = 0.1235476191121251
$ws.Range("P5").Value = 0.1235476191121251
$ws.Range("Q5").Value = 311.416715936375
$ws.Range("R5").Value = 2802.750443427375
$ws.Range("S5").Value = 0.01282043004400574
$ws.Range("T5").Value = 0.01282043004400574
# Row 6
$ws.Range("G6").Value = 47.94465366666667
$ws.Range("I6").Value = 0.5196887643218222
$ws.Range("J6").Value = 0.5196887643218222
$ws.Range("M6").Value = 227.11144
$ws.Range("N6").Value = 681.33432
$ws.Range("O6").Value = 0.8625743548356182
$ws.Range("P6").Value = 0.8625743548356182
$ws.Range("Q6").Value = 10888.77933453795
$ws.Range("R6").Value = 97999.01401084154
$ws.Range("S6").Value = 0.4482702006002154
$ws.Range("T6").Value = 0.4482702006002154
# Row 7
$ws.Range("G7").Value = 47.94465366666667
$ws.Range("I7").Value = 0.5196887643218222
$ws.Range("J7").Value = 0.5196887643218222
$ws.Range("O7").Value = 0.001598666154760757
$ws.Range("P7").Value = 0.001598666154760757
$ws.Range("S7").Value = 0.0008308088385307368
$ws.Range("T7").Value = 0.0008308088385307365
# Row 8
$ws.Range("G8").Value = 47.94465366666667
$ws.Range("I8").Value = 0.5196887643218222
$ws.Range("J8").Value = 0.5196887643218222
$ws.Range("M8").Value = 3.233093
$ws.Range("N8").Value = 9.699279000000001
$ws.Range("O8").Value = 0.01227935989749593
$ws.Range("P8").Value = 0.01227935989749593
$ws.Range("Q8").Value = 155.0095241571244
$ws.Range("R8").Value = 1395.085717414119
$ws.Range("S8").Value = 0.006381445371792598
$ws.Range("T8").Value = 0.006381445371792597
# Row 9
$ws.Range("G9").Value = 47.94465366666667
$ws.Range("I9").Value = 0.5196887643218222
$ws.Range("J9").Value = 0.5196887643218222
$ws.Range("M9").Value = 32.52945966666667
$ws.Range("N9").Value = 97.588379
$ws.Range("O9").Value = 0.1235476191121251
$ws.Range("P9").Value = 0.1235476191121251
$ws.Range("Q9").Value = 1559.613677682136
$ws.Range("R9").Value = 14036.52309913922
$ws.Range("S9").Value = 0.06420630951128346
$ws.Range("T9").Value = 0.06420630951128346
# Row 10
$ws.Range("G10").Value = 11.32006633333333
$ws.Range("H10").Value = 33.960199
$ws.Range("I10").Value = 0.122702133291269
$ws.Range("J10").Value = 0.122702133291269
$ws.Range("M10").Value = 227.11144
$ws.Range("N10").Value = 681.33432
$ws.Range("O10").Value = 0.8625743548356182
$ws.Range("P10").Value = 0.8625743548356182
$ws.Range("Q10").Value = 2570.916565858854
$ws.Range("R10").Value = 23138.24909272968
$ws.Range("S10").Value = 0.1058397134606704
$ws.Range("T10").Value = 0.1058397134606704
# Row 11
$ws.Range("G11").Value = 11.32006633333333
$ws.Range("H11").Value = 33.960199
$ws.Range("I11").Value = 0.122702133291269
$ws.Range("J11").Value = 0.122702133291269
$ws.Range("O11").Value = 0.001598666154760757
$ws.Range("P11").Value = 0.001598666154760757
$ws.Range("Q11").Value = 4.764849867737557
$ws.Range("R11").Value = 42.88364880963801
$ws.Range("S11").Value = 0.0001961597476096948
$ws.Range("T11").Value = 0.0001961597476096948
# Row 12
$ws.Range("G12").Value = 11.32006633333333
$ws.Range("H12").Value = 33.960199
$ws.Range("I12").Value = 0.122702133291269
$ws.Range("J12").Value = 0.122702133291269
$ws.Range("M12").Value = 3.233093
$ws.Range("N12").Value = 9.699279000000001
$ws.Range("O12").Value = 0.01227935989749593
$ws.Range("P12").Value = 0.01227935989749593
$ws.Range("Q12").Value = 36.59882722183568
$ws.Range("R12").Value = 329.389444996521
$ws.Range("S12").Value = 0.001506703654874009
$ws.Range("T12").Value = 0.001506703654874008
# Row 13
$ws.Range("G13").Value = 11.32006633333333
$ws.Range("H13").Value = 33.960199
$ws.Range("I13").Value = 0.122702133291269
$ws.Range("J13").Value = 0.122702133291269
$ws.Range("M13").Value = 32.52945966666667
$ws.Range("N13").Value = 97.588379
$ws.Range("O13").Value = 0.1235476191121251
$ws.Range("P13").Value = 0.1235476191121251
$ws.Range("Q13").Value = 368.235641214158
$ws.Range("R13").Value = 3314.120770927421
$ws.Range("S13").Value = 0.01515955642811491
$ws.Range("T13").Value = 0.01515955642811491
# Row 14
$ws.Range("G14").Value = 23.41838033333333
$ws.Range("H14").Value = 70.25514099999999
$ws.Range("I14").Value = 0.2538399635225604
$ws.Range("J14").Value = 0.2538399635225604
$ws.Range("M14").Value = 227.11144
$ws.Range("N14").Value = 681.33432
$ws.Range("O14").Value = 0.8625743548356182
$ws.Range("P14").Value = 0.8625743548356182
$ws.Range("Q14").Value = 5318.582079971014
$ws.Range("R14").Value = 47867.23871973912
$ws.Range("S14").Value = 0.2189558427669694
$ws.Range("T14").Value = 0.2189558427669694
# Row 15
$ws.Range("G15").Value = 23.41838033333333
$ws.Range("H15").Value = 70.25514099999999
$ws.Range("I15").Value = 0.2538399635225604
$ws.Range("J15").Value = 0.2538399635225604
$ws.Range("O15").Value = 0.001598666154760757
$ws.Range("P15").Value = 0.001598666154760757
$ws.Range("Q15").Value = 9.857280262160224
$ws.Range("R15").Value = 88.71552235944201
$ws.Range("S15").Value = 0.0004058053584092225
$ws.Range("T15").Value = 0.0004058053584092224
# Row 16
$ws.Range("G16").Value = 23.41838033333333
$ws.Range("H16").Value = 70.25514099999999
$ws.Range("I16").Value = 0.2538399635225604
$ws.Range("J16").Value = 0.2538399635225604
$ws.Range("M16").Value = 3.233093
$ws.Range("N16").Value = 9.699279000000001
$ws.Range("O16").Value = 0.01227935989749593
$ws.Range("P16").Value = 0.01227935989749593
$ws.Range("Q16").Value = 75.71380152703766
$ws.Range("R16").Value = 681.424213743339
$ws.Range("S16").Value = 0.003116992268460759
$ws.Range("T16").Value = 0.003116992268460758
# Row 17
$ws.Range("G17").Value = 23.41838033333333
$ws.Range("H17").Value = 70.25514099999999
$ws.Range("I17").Value = 0.2538399635225604
$ws.Range("J17").Value = 0.2538399635225604
$ws.Range("M17").Value = 32.52945966666667
$ws.Range("N17").Value = 97.588379
$ws.Range("O17").Value = 0.1235476191121251
$ws.Range("P17").Value = 0.1235476191121251
$ws.Range("Q17").Value = 761.7872585118265
$ws.Range("R17").Value = 6856.085326606439
$ws.Range("S17").Value = 0.03136132312872104
$ws.Range("T17").Value = 0.03136132312872104

Write-Output "Applied 178 cell updates"